$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.12
$ws.Range("G2").Value = 2.64
$ws.Range("H2").Value = 2.68
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 4.7
$ws.Range("L2").Value = 1.23
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 1.94
$ws.Range("Q2").Value = 1.71
$ws.Range("R2").Value = 1.36
$ws.Range("U2").Value = 1.89
$ws.Range("W2").Value = 1.6
$ws.Range("G3").Value = 1.62
$ws.Range("H3").Value = 5.3
$ws.Range("J3").Value = 4.7
$ws.Range("F4").Value = 2.72
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 2.46
$ws.Range("I4").Value = 2.84
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 4.5
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 1.5
$ws.Range("I6").Value = 1.52
$ws.Range("J6").Value = 4.7
$ws.Range("T6").Value = 1.81
$ws.Range("U6").Value = 2.08
$ws.Range("AA6").Value = 15.5
$ws.Range("AB6").Value = 32
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 12
$ws.Range("AE6").Value = 18
$ws.Range("AI6").Value = 36
$ws.Range("AJ6").Value = 250
$ws.Range("AO6").Value = 6.6
